# Update crypto price/volume data on Sheet1 per the latest GitHub Actions scrape.
# Values are written as literal text (matching the source inlineStr cells) rather
# than being auto-coerced to numbers/dates by Excel's input parser.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $r = $ws.Range($CellRef)
    # Force text interpretation so numeric-looking strings (e.g. "234.96",
    # "1.000", "30.317.58") and padded percentages are stored verbatim
    # instead of being parsed into numbers.
    $r.NumberFormat = "@"
    $r.Value = $Text
    # Restore the default "Normal" style so we don't leave a stray
    # text-format style applied to the cell (matches the unstyled source).
    $r.Style = "Normal"
}

Set-Text "D2" '30.317.58'
Set-Text "E2" '  +0.24%  '
Set-Text "D3" '1.865.66'
Set-Text "E3" '  +0.17%  '
Set-Text "D4" '1.000'
Set-Text "E4" '  -0.05%  '
Set-Text "D5" '234.96'
Set-Text "E5" '  -0.77%  '
Set-Text "D6" '1.000'
Set-Text "E6" '  -0.04%  '
Set-Text "D7" '0.4708'
Set-Text "E7" '  -0.35%  '
Set-Text "D8" '0.2862'
Set-Text "E8" '  -1.39%  '
Set-Text "D9" '0.06582'
Set-Text "E9" '  +0.48%  '
Set-Text "D10" '21.31'
Set-Text "E10" '  -2.48%  '
Set-Text "D11" '0.07836'
Set-Text "E11" '  -1.23%  '
Set-Text "D12" '97.10'
Set-Text "E12" '  -0.76%  '
Set-Text "D13" '1.862.69'
Set-Text "E13" '  -0.02%  '
Set-Text "D14" '0.6963'
Set-Text "E14" '  +2.06%  '
Set-Text "D15" '5.099'
Set-Text "D16" '269.48'
Set-Text "E16" '  +0.85%  '
Set-Text "D17" '30.279.16'
Set-Text "E17" '  +0.13%  '
Set-Text "D18" '13.85'
Set-Text "E18" '  +0.73%  '
Set-Text "D19" '0.000007645'
Set-Text "E19" '  +3.24%  '
Set-Text "D20" '1.001'
Set-Text "E20" '  +0.07%  '
Set-Text "D21" '2.112.99'
Set-Text "E21" '  +0.13%  '
Set-Text "E22" '  -0.02%  '
Set-Text "D23" '5.232'
Set-Text "E23" '  -1.52%  '
Set-Text "D24" '6.158'
Set-Text "E24" '  -0.47%  '
Set-Text "D25" '9.446'
Set-Text "E25" '  +2.33%  '
Set-Text "D26" '167.27'
Set-Text "E26" '  -0.18%  '
Set-Text "D27" '18.90'
Set-Text "E27" '  -0.14%  '
Set-Text "D28" '1.942'
Set-Text "E28" '  -1.20%  '
Set-Text "B29" 'Stellar'
Set-Text "C29" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-Text "D29" '0.09947'
Set-Text "E29" '  +0.83%  '
Set-Text "B30" 'Toncoin'
Set-Text "C30" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-Text "D30" '1.363'
Set-Text "E30" '  -2.24%  '
Set-Text "D31" '4.357'
Set-Text "E31" '  -0.42%  '
Set-Text "D32" '1.459'
Set-Text "E32" '  -0.87%  '
Set-Text "D33" '4.056'
Set-Text "E33" '  -0.08%  '
Set-Text "D34" '0.04730'
Set-Text "E34" '  +0.28%  '
Set-Text "E35" '  +0.54%  '
Set-Text "D36" '0.7022'
Set-Text "E36" '  -0.28%  '
Set-Text "D37" '2.711'
Set-Text "E37" '  +0.16%  '
Set-Text "E38" '  -0.36%  '
Set-Text "D39" '2.791'
Set-Text "E39" '  +7.02%  '
Set-Text "D40" '6.315'
Set-Text "E40" '  +0.84%  '
Set-Text "D41" '73.05'
Set-Text "E41" '  -1.71%  '
Set-Text "D42" '1.951'
Set-Text "E42" '  +0.47%  '
Set-Text "D43" '0.4174'
Set-Text "E43" '  +0.08%  '
Set-Text "E44" '  +0.05%  '
Set-Text "D45" '0.8368'
Set-Text "E45" '  -1.01%  '
Set-Text "D46" '103.20'
Set-Text "E46" '  -0.31%  '
Set-Text "D47" '970.56'
Set-Text "E47" '  +2.01%  '
Set-Text "D48" '7.123'
Set-Text "E48" '  -0.79%  '
Set-Text "D49" '9.161'
Set-Text "E49" '  -1.06%  '
Set-Text "D50" '34.48'
Set-Text "E50" '  +0.95%  '
Set-Text "E51" '  +0.38%  '
